$wb = $excel.ActiveWorkbook

# Sheet "2025" (first sheet) - row 2 updates
$ws1 = $wb.Worksheets.Item("2025")
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = 0.009216388360670347
$ws1.Range("E2").Value = 0.3881565890751533
$ws1.Range("G2").Value = 0.2494892361374987
$ws1.Range("I2").Value = 0.3866351
$ws1.Range("L2").Value = 0.5980458758222187
$ws1.Range("M2").Value = 0.08407824999999999
$ws1.Range("N2").Value = 10.84655480891572
$ws1.Range("O2").Value = 2.862907308505379

# Sheet "2030" (second sheet) - row 2 updates
$ws2 = $wb.Worksheets.Item("2030")
$ws2.Range("A2").Value = 0.05356848075857206
$ws2.Range("B2").Value = 0.06920837885335457
$ws2.Range("E2").Value = 0.4138138
$ws2.Range("I2").Value = 0.8506273458905896
$ws2.Range("L2").Value = 0.1483279332818203
$ws2.Range("M2").Value = 0.1014180740827377
$ws2.Range("N2").Value = 15.13001910209066
$ws2.Range("O2").Value = 7.28386422727996
